$d = $word.ActiveDocument

# Locate the phrase containing the typo ("mob up" should read "mop up").
$found = $d.Content
$found.Find.Execute("can mob up", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)

# The 'b' in "mob" is the 7th character of "can mob up" (index 6).
$bStart = $found.Start + 6
$bEnd = $bStart + 1

# Replace the single mistyped character, splitting the run at that point.
$charRange = $d.Range($bStart, $bEnd)
$charRange.Text = "p"

# Force the newly typed character into its own run (matching how Word
# marks the most-recently-edited text) while keeping formatting identical
# to its neighbours.
$newRunRange = $d.Range($bStart, $bStart + 1)
$newRunRange.Font.Size = $newRunRange.Font.Size + 1
$newRunRange.Font.Size = $newRunRange.Font.Size - 1

# Word automatically keeps a single hidden "_GoBack" bookmark marking the
# location of the last edit. Move it from its old location to right after
# the corrected character.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$goBackRange = $d.Range($bStart + 1, $bStart + 1)
$d.Bookmarks.Add("_GoBack", $goBackRange)
